# Fix a misspelled word: "serían" (plural) -> "sería" (singular) in the
# paragraph discussing microservices / architecture.
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("serían", $true, $false, $false, $false, $false, $true, 1, $false, "sería", 2)
